$d = $word.ActiveDocument
$d.Content.Find.Execute("4+ years of progressive", $true, $false, $false, $false, $false, $true, 1, $false, "5+ years of progressive", 2)
